# fix kbv report table
# The "Regional Geimpfte" sheet had two stray/mislabeled columns ("cases" and
# "R0") in between the percentage columns and the actual 7-Tage-Inzidenz
# columns. Remove those two columns so that the former H/I ("7-Tage-Inzidenz"
# / "7-Tage-Inzidenz 60+") columns shift left into F/G.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Regional Geimpfte")
$ws.Range("F1:G1").EntireColumn.Delete()
